# Append the "Build Fix Update" report block to the end of the document,
# mirroring the existing report blocks already present (same run/paragraph
# formatting: Helvetica Light, sz 24).
#
# We build the new content as a literal OOXML fragment (so tab stops stay
# real <w:tab/> runs rather than literal tab characters, and empty text
# runs keep xml:space="preserve") and splice it in right after the very
# last paragraph via Range.InsertXML.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr>'

$p1 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + '<w:t xml:space="preserve"></w:t></w:r></w:p>'
$p2 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + '<w:t xml:space="preserve">---</w:t></w:r></w:p>'
$p3 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + '<w:t xml:space="preserve">Build Fix Update</w:t></w:r></w:p>'
$p4 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + '<w:t xml:space="preserve">Updated: 2026-02-18</w:t></w:r></w:p>'
$p5 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + '<w:t xml:space="preserve"></w:t></w:r></w:p>'
$p6 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + `
        '<w:t xml:space="preserve">Module Name</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">Developed</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">Partial Developed</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">Need To Develop</w:t>' + `
        '</w:r></w:p>'
$p7 = '<w:p ' + $wNs + '><w:pPr/><w:r>' + $rPr + `
        '<w:t xml:space="preserve">Migration Build Stability</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">Fixed storefront CMS migration compile syntax for CI</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">-</w:t><w:tab/>' + `
        '<w:t xml:space="preserve">-</w:t>' + `
        '</w:r></w:p>'

$fragment = $p1 + $p2 + $p3 + $p4 + $p5 + $p6 + $p7

# Use a free-standing zero-length Range at the absolute end of the story
# (NOT Paragraph.Range collapsed to its end) -- a collapsed Paragraph.Range
# stays associated with that paragraph and InsertXML there clobbers it;
# Document.Range(end, end) is a plain insertion point and appends cleanly.
$endPos = $d.Content.End
$tail = $d.Range($endPos, $endPos)
$null = $tail.InsertXML($fragment)
